$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Refresh the cryptos list snapshot (price + 1h volume change columns).
# A couple of rows (39/40) also swapped coin identity (Coin name + Link)
# along with their price/volume, matching the source feed's new ordering.
# For Price values that look like plain numbers, force the cell to Text
# first (then restore the Normal style) so Excel doesn't auto-convert the
# typed string into a floating-point number / scientific notation.

$ws.Range("D2").Value2 = '26.179.05'
$ws.Range("E2").Value2 = '  -0.56%  '
$ws.Range("D3").Value2 = '1.584.95'
$ws.Range("E3").Value2 = '  -0.39%  '
$ws.Range("E4").Value2 = '  -0.16%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value2 = '211.58'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value2 = '  +0.82%  '
$ws.Range("E6").Value2 = '  -0.05%  '
$ws.Range("E7").Value2 = '  -0.15%  '
$ws.Range("E8").Value2 = '  -0.43%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value2 = '0.0603'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value2 = '  -1.29%  '
$ws.Range("E10").Value2 = '  -1.95%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value2 = '0.0847'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value2 = '  +0.32%  '
$ws.Range("D12").Value2 = '1.807.66'
$ws.Range("E12").Value2 = '  -0.46%  '
$ws.Range("D13").Value2 = '1.596.11'
$ws.Range("E13").Value2 = '  +0.13%  '
$ws.Range("E14").Value2 = '  -1.78%  '
$ws.Range("E15").Value2 = '  -0.24%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value2 = '63.86'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value2 = '  -1.03%  '
$ws.Range("D17").Value2 = '26.173.28'
$ws.Range("E17").Value2 = '  -0.67%  '
$ws.Range("D18").Value2 = '0.0₃0723'
$ws.Range("E18").Value2 = '  -0.64%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value2 = '213.82'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value2 = '  +1.28%  '
$ws.Range("E20").Value2 = '  -2.78%  '
$ws.Range("E21").Value2 = '  -0.09%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value2 = '4.24'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value2 = '  -0.61%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value2 = '8.95'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value2 = '  +0.31%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value2 = '2.12'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value2 = '  -2.16%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value2 = '144.26'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value2 = '  -0.63%  '
$ws.Range("E26").Value2 = '  -0.13%  '
$ws.Range("E27").Value2 = '  -1.18%  '
$ws.Range("E28").Value2 = '  -0.94%  '
$ws.Range("E29").Value2 = '  -1.36%  '
$ws.Range("E30").Value2 = '  -2.07%  '
$ws.Range("E31").Value2 = '  +0.59%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value2 = '3.17'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value2 = '  -1.09%  '
$ws.Range("D33").Value2 = '1.408.22'
$ws.Range("E33").Value2 = '  +7.85%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value2 = '2.94'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value2 = '  -1.83%  '
$ws.Range("E35").Value2 = '  -0.69%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value2 = '0.587'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value2 = '  -4.37%  '
$ws.Range("E37").Value2 = '  -1.42%  '
$ws.Range("E38").Value2 = '  -1.51%  '
$ws.Range("B39").Value2 = 'FraxShare'
$ws.Range("C39").Value2 = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value2 = '5.87'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value2 = '  +4.35%  '
$ws.Range("B40").Value2 = 'ARBITRUM'
$ws.Range("C40").Value2 = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value2 = '0.818'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value2 = '  +0.83%  '
$ws.Range("E41").Value2 = '  -0.14%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value2 = '0.940'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value2 = '  -14.95%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value2 = '0.765'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value2 = '  +0.15%  '
$ws.Range("E44").Value2 = '  -0.31%  '
$ws.Range("D45").Value2 = '1.719.56'
$ws.Range("E45").Value2 = '  -0.53%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value2 = '60.86'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value2 = '  -2.79%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value2 = '85.16'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value2 = '  -3.21%  '
$ws.Range("E48").Value2 = '  -1.11%  '
$ws.Range("E49").Value2 = '  -0.87%  '
$ws.Range("E50").Value2 = '  -1.34%  '
$ws.Range("E51").Value2 = '  -0.25%  '
